$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from existing header cell (H1) to the new header cells
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Data for columns I (I0) and J (IF), keyed by row number
$data = @{
    2  = @(1, 3)
    3  = @(1, 5)
    4  = @(1, 4)
    5  = @(1, 6)
    6  = @(1, 6)
    7  = @(1, 6)
    8  = @(1, 6)
    9  = @(1, 6)
    10 = @(1, 5)
    11 = @(1, 5)
    12 = @(1, 5)
    13 = @(4, 5)
    14 = @(1, 6)
    15 = @(1, 5)
    16 = @(1, 4)
    17 = @(1, 3)
    18 = @(4, 5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}

$wb.Save()
